# metrica 2 dati ok
# Adds a "Media" (average) column L (row averages of B:K) and a "Media" summary
# row 24 (L24 label, M24 = average of the L column) to each of the 10 sheets.

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le 10; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Column L: per-row average of B:K for rows 1..20 (shared formula style,
    # exactly like Excel's "fill down" behaviour when you enter a formula in
    # L1 and then copy/fill it through L2:L20).
    $ws.Range("L1").Formula = "=AVERAGE(B1:K1)"
    $ws.Range("L2:L20").Formula = "=AVERAGE(B2:K2)"

    # Row 24: a "Media" label in L24 and the overall average of L1:L20 in M24.
    $ws.Range("L24").Value = "Media"
    $ws.Range("M24").Formula = "=AVERAGE(L1:L20)"
}

# Restore selections: every sheet now has the new L1:M24 block selected,
# except the last-touched sheet (Foglio9), whose cursor ends up on L28, and
# Foglio1, which becomes the active tab (previously Foglio10 was active).
for ($i = 2; $i -le 8; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Activate() | Out-Null
    $ws.Range("L1:M24").Select() | Out-Null
}

$ws10 = $wb.Worksheets.Item(10)
$ws10.Activate() | Out-Null
$ws10.Range("L1:M24").Select() | Out-Null

$ws9 = $wb.Worksheets.Item(9)
$ws9.Activate() | Out-Null
$ws9.Range("L28").Select() | Out-Null

$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate() | Out-Null
$ws1.Range("L1:M24").Select() | Out-Null
